# Generate Report for Handoff
#
# Updates the localization-status report so that the second tracked file
# (c5ebdeb9-67c7-4d04-8253-e671d9c6ce89.md) shows as "Ready for handoff"
# with refreshed handoff timestamps, and records an Error Detail noting
# that the handback file on record is stale (a newer commit exists).

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$handoffDateTimeOverview = "2016-09-01 04:52:48"
$handoffDateTimeZhCn = "2016-09-01 04:52:43"
$handoffDateTimeDeDe = "2016-09-01 04:52:48"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/05c210319cd32e806ca946c084a3efa24aa307ed/e2e/c5ebdeb9-67c7-4d04-8253-e671d9c6ce89.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a582f2b3270d5d438df800f29ecd83ec1ac4977/e2e/c5ebdeb9-67c7-4d04-8253-e671d9c6ce89.md."

# --- Overview sheet: row 3 is the c5ebdeb9 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $handoffDateTimeOverview

# Excel's ColumnWidth property (character units) is stored internally with a
# constant ~5/6 character padding offset, so assigning 40 directly serialises
# as width="40.833333333333336". Back the offset out so the saved OOXML
# <col> width attribute lands on exactly 40, matching the target column.
$colWidth40 = 40 - (5/6)

# --- zh-cn sheet: row 3 is the c5ebdeb9 file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("H3").Value = $handoffDateTimeZhCn
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = $colWidth40

# --- de-de sheet: row 3 is the c5ebdeb9 file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("H3").Value = $handoffDateTimeDeDe
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $colWidth40
